$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 66 (Serie 01-01-2021): Pasivos netos incurridos / Endeudamiento Externo Neto
# changed, which flows through to the totals (AA) and difference (AB) columns.
$ws.Range("P66").Value = 7013
$ws.Range("AA66").Value = 7013
$ws.Range("AB66").Value = -7013

# Update row 67 (Serie 01-04-2021) similarly.
$ws.Range("P67").Value = 6444
$ws.Range("AA67").Value = 6444
$ws.Range("AB67").Value = -6444

# Append new row 68 for the new quarterly series 01-07-2021.
# Force column A to be stored as text (not auto-converted to a date serial),
# matching the string cells used for every other "Serie" column value, then
# clear the temporary number format so no extra cell style is introduced.
$ws.Range("A68").NumberFormat = "@"
$ws.Range("A68").Value = "01-07-2021"
$ws.Range("A68").ClearFormats()

$ws.Range("B68").Value = 0
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("N68").Value = 0
$ws.Range("O68").Value = 0
$ws.Range("P68").Value = 5866
$ws.Range("Q68").Value = 0
$ws.Range("R68").Value = 0
$ws.Range("S68").Value = 0
$ws.Range("T68").Value = 0
$ws.Range("U68").Value = 0
$ws.Range("V68").Value = 0
$ws.Range("W68").Value = 0
$ws.Range("X68").Value = 0
$ws.Range("Y68").Value = 0
$ws.Range("Z68").Value = 0
$ws.Range("AA68").Value = 5866
$ws.Range("AB68").Value = -5866
